$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, for columns E,G,H,K,M,N,O,P,Q,R,S,T
$data = @{
    2 = @{ E=3; G=1.607064; H=4.821192; K=3; M=11.88831466666667; N=35.664944; O=0.1647867559347935; P=0.1715303447251208; Q=19.105282521472; R=171.947542693248; S=0.1647867559347935; T=0.1715303447251208 }
    3 = @{ E=3; G=1.607064; H=4.821192; K=3; M=17.52633466666667; N=52.579004; O=0.2429366915434531; P=0.2528784198125617; Q=28.165941494752; R=253.493473452768; S=0.2429366915434531; T=0.2528784198125617 }
    4 = @{ E=3; G=1.607064; H=4.821192; K=3; M=13.957045; N=41.871135; O=0.1934619189071989; P=0.2013789849377604; Q=22.42986456588; R=201.86878109292; S=0.1934619189071989; T=0.2013789849377604 }
    5 = @{ E=3; G=1.607064; H=4.821192; K=3; M=20.26311466666667; N=60.789344; O=0.2808718497683384; P=0.2923660032084714; Q=32.564122108672; R=293.077098978048; S=0.2808718497683384; T=0.2923660032084714 }
    6 = @{ E=3; G=1.607064; H=4.821192; K=2; M=8.508820500000001; N=17.017641; O=0.1179427838462161; P=0.08184624731608578; Q=13.674219108012; R=82.04531464807201; S=0.1179427838462161; T=0.08184624731608578 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

$wb.Save()
